$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at S:T (shifts old S..AC to U..AE), creating room
# for the new "v.background"/"v.animate"/"v.mc" climbview columns.
$ws.Range("S:T").EntireColumn.Insert()

# Rename the existing "visual" header to "v.background" (same column, new label).
$ws.Range("N1").Value = "v.background"

# New header labels for the two freshly inserted columns.
$ws.Range("S1").Value = "v.animate"
$ws.Range("T1").Value = "v.mc"

# New data point introduced alongside the climbview columns.
$ws.Range("T3").Value = 1

# Row 5 keeps its old mc1: monitor value in S5 *and* gets the same value
# duplicated into the new v.mc column at U5.
$ws.Range("S5").Value = "https://i.imgsafe.org/c5fd4eaf6b.png"
$ws.Range("U5").Value = "https://i.imgsafe.org/c5fd4eaf6b.png"

# Selection moves to T3, matching the saved worksheet view.
$ws.Range("T3").Select() | Out-Null
